$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52: move "Scalpel Accuracy:" label from C52 to E52, clear D52,
# and set new accuracy value in F52 (was empty).
$ws.Range("C52").Value = $null
$ws.Range("D52").Value = $null
$ws.Range("E52").Value = "Scalpel Accuracy:"
$ws.Range("F52").Value = 91.84

# Row 53: rename label text, value in F53 stays the same.
$ws.Range("E53").Value = "Accuracy vs PyType"
